$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert a new "2022-Q4" sheet right after "总计", by copying the existing
#    "2022-Q3" sheet (same column layout/styling) and then overwriting its
#    data with the 2022-Q4 figures.
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$total = $wb.Worksheets.Item("总计")
$q3.Copy($null, $total)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# 2022-Q3 had two funds; 2022-Q4 only has one (457001), so drop the second
# fund's row.
$q4.Rows.Item(3).Delete()

# Overwrite the remaining fund's figures with the 2022-Q4 values (these
# columns are stored as text in the workbook, hence the leading "'").
$q4.Range("D2").Value = "'3.18"
$q4.Range("E2").Value = "'88.46"
$q4.Range("F2").Value = "'2.96"
$q4.Range("G2").Value = "'0.0941"
$q4.Range("H2").Value = 10

# Re-apply the (unstyled) number formatting from C2 so the text cells don't
# pick up a stray "quote prefix" style.
$q4.Range("C2").Copy() | Out-Null
$q4.Range("D2:G2").PasteSpecial(-4122) | Out-Null
$q4.Range("A1").Select() | Out-Null

# ---------------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: insert a new row for 2022-Q4 right
#    after the header, and renumber the existing rows' index column.
# ---------------------------------------------------------------------------
$total.Rows.Item(2).Insert()
$total.Range("A2:D2").ClearFormats()
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.09

# Match formatting (A has the bordered/centred style, B:D stay plain) from
# the row below, which already has the correct look.
$total.Range("A3").Copy() | Out-Null
$total.Range("A2").PasteSpecial(-4122) | Out-Null

# The other rows shift down one place; keep their data but renumber the
# zero-based index in column A.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5
$total.Range("A1").Select() | Out-Null

# ---------------------------------------------------------------------------
# 3) Restore the originally-active tab ("2020-Q4", now the last sheet).
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("2020-Q4").Activate() | Out-Null
